# Applies the cryptocurrency price/volume data refresh described by the commit
# "Updated cryptos list on Tue Jun  4 08:13:58 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price-column value while keeping it stored as plain text
# (mirrors the source data which stores prices as text, e.g. "6.80"), and
# without leaving any residual explicit cell style/number-format behind.
function Set-PriceText($rowNum, $text) {
    $cell = $ws.Cells.Item($rowNum, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-VolumeText($rowNum, $text) {
    $ws.Cells.Item($rowNum, 5).Value = $text
}

# --- Per-row updates (Price + Volume(1h)) ---
Set-PriceText 2 "68.966.15"
Set-VolumeText 2 "  -0.02%  "

Set-PriceText 3 "3.770.28"
Set-VolumeText 3 "  -1.29%  "

Set-VolumeText 4 "  -0.01%  "

Set-PriceText 5 "629.56"
Set-VolumeText 5 "  +0.32%  "

Set-PriceText 6 "165.38"
Set-VolumeText 6 "  -0.03%  "

Set-PriceText 7 "3.768.74"
Set-VolumeText 7 "  -1.26%  "

Set-VolumeText 8 "  +0.06%  "

Set-VolumeText 9 "  -0.31%  "

Set-VolumeText 10 "  -2.02%  "

Set-VolumeText 11 "  +0.40%  "

Set-PriceText 12 "6.80"
Set-VolumeText 12 "  +2.81%  "

Set-VolumeText 13 "  -4.66%  "

Set-PriceText 14 "34.78"
Set-VolumeText 14 "  -3.25%  "

Set-PriceText 15 "4.403.72"
Set-VolumeText 15 "  -1.23%  "

Set-PriceText 16 "3.780.83"
Set-VolumeText 16 "  -1.60%  "

Set-PriceText 17 "68.944.09"
Set-VolumeText 17 "  -0.02%  "

Set-PriceText 18 "17.63"
Set-VolumeText 18 "  -2.66%  "

Set-VolumeText 19 "  -0.20%  "

Set-VolumeText 20 "  -1.94%  "

Set-PriceText 21 "467.50"
Set-VolumeText 21 "  +0.48%  "

Set-VolumeText 22 "  -1.81%  "

Set-PriceText 23 "0.701"
Set-VolumeText 23 "  -0.99%  "

Set-PriceText 24 "81.96"
Set-VolumeText 24 "  -2.25%  "

Set-PriceText 25 "0.0000142"
Set-VolumeText 25 "  -7.90%  "

Set-PriceText 26 "12.12"
Set-VolumeText 26 "  +0.81%  "

Set-VolumeText 27 "  -1.54%  "

Set-PriceText 28 "10.11"
Set-VolumeText 28 "  +0.54%  "

Set-VolumeText 29 "  -0.03%  "

Set-PriceText 30 "3.919.33"
Set-VolumeText 30 "  -1.21%  "

Set-PriceText 31 "2.28"
Set-VolumeText 31 "  +2.10%  "

Set-VolumeText 32 "  -0.43%  "

Set-PriceText 33 "7.10"
Set-VolumeText 33 "  -2.91%  "

Set-PriceText 34 "0.178"
Set-VolumeText 34 "  +19.84%  "

Set-PriceText 35 "28.43"
Set-VolumeText 35 "  -2.63%  "

Set-PriceText 36 "1.00"
Set-VolumeText 36 "  -0.12%  "

Set-PriceText 37 "3.723.47"
Set-VolumeText 37 "  -1.10%  "

Set-PriceText 38 "8.87"
Set-VolumeText 38 "  -2.30%  "

Set-PriceText 39 "0.101"
Set-VolumeText 39 "  -0.62%  "

Set-VolumeText 40 "  -4.32%  "

Set-VolumeText 41 "  -2.41%  "

Set-VolumeText 45 "  +5.68%  "

Set-PriceText 46 "156.14"
Set-VolumeText 46 "  +1.17%  "

Set-PriceText 47 "43.62"
Set-VolumeText 47 "  +1.82%  "

Set-PriceText 48 "46.99"
Set-VolumeText 48 "  +0.24%  "

Set-VolumeText 49 "  -2.54%  "

Set-VolumeText 50 "  -2.34%  "

Set-PriceText 51 "8.35"
Set-VolumeText 51 "  -1.25%  "

# Rows 42 and 43 swapped order (Mantle <-> FirstDigitalUSD) along with updated data
$ws.Cells.Item(42, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-PriceText 42 "1.00"
Set-VolumeText 42 "  +0.02%  "

$ws.Cells.Item(43, 2).Value = "Mantle"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText 43 "0.959"
Set-VolumeText 43 "  -2.35%  "

